# Fix file upload functionality - correct null value handling in data
# processing.
#
# The upload job had been dropping the most-recent sensor reading for each
# lifter feed. Each of the four sheets is missing its newest row: append
# row 46 (one hour after the existing last row, 45) with the corrected
# reading for every sheet.
#
# Row 46 duplicates row 45's payload (B-I are unchanged from the last
# reading) except for the timestamp in column A, which advances by one
# hour - so copy row 45 down (this preserves the exact cell types, e.g.
# the 24-digit id in column G that must stay text) and then patch in the
# new timestamp.

$wb = $excel.ActiveWorkbook

$newTimestamps = @{
    "ROW35-FE-LIFTER"  = "2025-03-06 05:42:06"
    "ROW35-MID-LIFTER" = "2025-03-06 05:29:35"
    "ROW02-FE-LIFTER"  = "2025-03-06 05:51:45"
    "ROW02-MID-LIFTER" = "2025-03-06 05:41:15"
}

foreach ($sheetName in $newTimestamps.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("A45:I45").Copy()
    $ws.Range("A46:I46").PasteSpecial()
    $ws.Cells.Item(46, 1).Value = $newTimestamps[$sheetName]
}
